$wb = $excel.ActiveWorkbook

$ws_ALC = $wb.Worksheets.Item("ALC")
$ws_ARM = $wb.Worksheets.Item("ARM")
$ws_BSM = $wb.Worksheets.Item("BSM")
$ws_CRP = $wb.Worksheets.Item("CRP")
$ws_CUL = $wb.Worksheets.Item("CUL")
$ws_GSM = $wb.Worksheets.Item("GSM")
$ws_LTW = $wb.Worksheets.Item("LTW")
$ws_WVR = $wb.Worksheets.Item("WVR")

$ws_ALC.Range("H129").Value = 1070.5652
$ws_ALC.Range("J129").Value = 1105.164
$ws_ALC.Range("L129").Value = 3315.492
$ws_ALC.Range("N129").Value = -13315.492
$ws_ALC.Range("H135").Value = 2560.4375
$ws_ALC.Range("I135").Value = 1984.4138
$ws_ALC.Range("J135").Value = 8128.6665
$ws_ALC.Range("K135").Value = 17859.7242
$ws_ALC.Range("L135").Value = 73157.9985
$ws_ALC.Range("M135").Value = -15324.7242
$ws_ALC.Range("N135").Value = -78227.9985
$ws_ALC.Range("H137").Value = 1401.8269
$ws_ALC.Range("I137").Value = 1167.186
$ws_ALC.Range("J137").Value = 2522.889
$ws_ALC.Range("K137").Value = 3501.558
$ws_ALC.Range("L137").Value = 7568.667
$ws_ALC.Range("M137").Value = -951.558
$ws_ALC.Range("N137").Value = -12668.667
$ws_ALC.Range("H138").Value = 948.67
$ws_ALC.Range("I138").Value = 466.9
$ws_ALC.Range("J138").Value = 1430.44
$ws_ALC.Range("K138").Value = 1400.7
$ws_ALC.Range("L138").Value = 4291.32
$ws_ALC.Range("N138").Value = -14571.32
$ws_ALC.Range("H141").Value = 1468
$ws_ALC.Range("I141").Value = 1010.3158
$ws_ALC.Range("J141").Value = 2917.3333
$ws_ALC.Range("K141").Value = 3030.9474
$ws_ALC.Range("L141").Value = 8751.999899999999
$ws_ALC.Range("M141").Value = 2149.0526
$ws_ALC.Range("N141").Value = -19111.9999
$ws_ARM.Range("H26").Value = 5201.1665
$ws_ARM.Range("I26").Value = 3801.75
$ws_ARM.Range("K26").Value = 3801.75
$ws_ARM.Range("M26").Value = -3471.75
$ws_ARM.Range("H32").Value = 4071.3
$ws_ARM.Range("I32").Value = 3077.1765
$ws_ARM.Range("J32").Value = 9704.666999999999
$ws_ARM.Range("K32").Value = 3077.1765
$ws_ARM.Range("L32").Value = 9704.666999999999
$ws_ARM.Range("M32").Value = -2790.1765
$ws_ARM.Range("N32").Value = -10278.667
$ws_ARM.Range("H61").Value = 4509.943
$ws_ARM.Range("I61").Value = 4617.3335
$ws_ARM.Range("K61").Value = 4617.3335
$ws_ARM.Range("M61").Value = -4405.3335
$ws_ARM.Range("H74").Value = 1056.8591
$ws_ARM.Range("I74").Value = 801.9167
$ws_ARM.Range("K74").Value = 801.9167
$ws_ARM.Range("M74").Value = 72.08330000000001
$ws_ARM.Range("H77").Value = 1056.8591
$ws_ARM.Range("I77").Value = 801.9167
$ws_ARM.Range("K77").Value = 4009.5835
$ws_ARM.Range("M77").Value = 358.4165000000003
$ws_ARM.Range("H88").Value = 2349.125
$ws_ARM.Range("I88").Value = 1974
$ws_ARM.Range("J88").Value = 2724.25
$ws_ARM.Range("K88").Value = 1974
$ws_ARM.Range("L88").Value = 2724.25
$ws_ARM.Range("M88").Value = -1568
$ws_ARM.Range("N88").Value = -3536.25
$ws_ARM.Range("H91").Value = 2349.125
$ws_ARM.Range("I91").Value = 1974
$ws_ARM.Range("J91").Value = 2724.25
$ws_ARM.Range("K91").Value = 1974
$ws_ARM.Range("L91").Value = 2724.25
$ws_ARM.Range("M91").Value = -570
$ws_ARM.Range("N91").Value = -5532.25
$ws_ARM.Range("H136").Value = 4509.943
$ws_ARM.Range("I136").Value = 4617.3335
$ws_ARM.Range("K136").Value = 13852.0005
$ws_ARM.Range("M136").Value = -11302.0005
$ws_BSM.Range("H20").Value = 9037.807000000001
$ws_BSM.Range("I20").Value = 808.8570999999999
$ws_BSM.Range("J20").Value = 15814.588
$ws_BSM.Range("K20").Value = 808.8570999999999
$ws_BSM.Range("L20").Value = 15814.588
$ws_BSM.Range("M20").Value = -561.8570999999999
$ws_BSM.Range("N20").Value = -16308.588
$ws_BSM.Range("H134").Value = 3045.4546
$ws_BSM.Range("I134").Value = 3187.4695
$ws_BSM.Range("K134").Value = 9562.408500000001
$ws_BSM.Range("M134").Value = -7027.408500000001
$ws_CRP.Range("H35").Value = 1951.125
$ws_CRP.Range("I35").Value = 1951.125
$ws_CRP.Range("K35").Value = 1951.125
$ws_CRP.Range("M35").Value = -1657.125
$ws_CRP.Range("H99").Value = 25038802
$ws_CRP.Range("I99").Value = 54670.668
$ws_CRP.Range("J99").Value = 62515000
$ws_CRP.Range("K99").Value = 54670.668
$ws_CRP.Range("L99").Value = 62515000
$ws_CRP.Range("M99").Value = -53172.668
$ws_CRP.Range("N99").Value = -62517996
$ws_CRP.Range("H126").Value = 25038802
$ws_CRP.Range("I126").Value = 54670.668
$ws_CRP.Range("J126").Value = 62515000
$ws_CRP.Range("K126").Value = 164012.004
$ws_CRP.Range("L126").Value = 187545000
$ws_CRP.Range("M126").Value = -161542.004
$ws_CRP.Range("N126").Value = -187549940
$ws_CRP.Range("H132").Value = 1989.2424
$ws_CRP.Range("I132").Value = 1479.4117
$ws_CRP.Range("J132").Value = 3722.6667
$ws_CRP.Range("K132").Value = 4438.2351
$ws_CRP.Range("L132").Value = 11168.0001
$ws_CRP.Range("M132").Value = -1908.2351
$ws_CRP.Range("N132").Value = -16228.0001
$ws_CRP.Range("H134").Value = 1960.0702
$ws_CRP.Range("I134").Value = 2193.3684
$ws_CRP.Range("K134").Value = 6580.1052
$ws_CRP.Range("M134").Value = -4045.1052
$ws_CUL.Range("H22").Value = 2666.6667
$ws_CUL.Range("J22").Value = 3040
$ws_CUL.Range("L22").Value = 9120
$ws_CUL.Range("N22").Value = -9458
$ws_CUL.Range("H26").Value = 31250104
$ws_CUL.Range("J26").Value = 38461628
$ws_CUL.Range("L26").Value = 115384884
$ws_CUL.Range("N26").Value = -115385460
$ws_CUL.Range("H27").Value = 2666.6667
$ws_CUL.Range("J27").Value = 3040
$ws_CUL.Range("L27").Value = 9120
$ws_CUL.Range("N27").Value = -9324
$ws_CUL.Range("H32").Value = 4400
$ws_CUL.Range("J32").Value = 4400
$ws_CUL.Range("L32").Value = 13200
$ws_CUL.Range("N32").Value = -13766
$ws_CUL.Range("H34").Value = 1397.6666
$ws_CUL.Range("I34").Value = 936
$ws_CUL.Range("J34").Value = 1974.75
$ws_CUL.Range("K34").Value = 2808
$ws_CUL.Range("L34").Value = 5924.25
$ws_CUL.Range("M34").Value = -2724
$ws_CUL.Range("N34").Value = -6092.25
$ws_CUL.Range("H38").Value = 4545600
$ws_CUL.Range("I38").Value = 7142879.5
$ws_CUL.Range("J38").Value = 361.5
$ws_CUL.Range("K38").Value = 21428638.5
$ws_CUL.Range("L38").Value = 1084.5
$ws_CUL.Range("M38").Value = -21428291.5
$ws_CUL.Range("N38").Value = -1778.5
$ws_CUL.Range("H39").Value = 1900
$ws_CUL.Range("J39").Value = 1900
$ws_CUL.Range("L39").Value = 5700
$ws_CUL.Range("N39").Value = -6288
$ws_CUL.Range("H40").Value = 252.80952
$ws_CUL.Range("I40").Value = 128.27777
$ws_CUL.Range("J40").Value = 1000
$ws_CUL.Range("K40").Value = 513.11108
$ws_CUL.Range("L40").Value = 4000
$ws_CUL.Range("M40").Value = -444.11108
$ws_CUL.Range("N40").Value = -4138
$ws_CUL.Range("H46").Value = 2500
$ws_CUL.Range("J46").Value = 2500
$ws_CUL.Range("L46").Value = 7500
$ws_CUL.Range("N46").Value = -7682
$ws_CUL.Range("H92").Value = 602
$ws_CUL.Range("J92").Value = 602
$ws_CUL.Range("L92").Value = 1806
$ws_CUL.Range("N92").Value = -4302
$ws_CUL.Range("H93").Value = 0
$ws_CUL.Range("J93").Value = 0
$ws_CUL.Range("L93").Value = 0
$ws_CUL.Range("H94").Value = 2982.8572
$ws_CUL.Range("I94").Value = 1940
$ws_CUL.Range("J94").Value = 3400
$ws_CUL.Range("K94").Value = 5820
$ws_CUL.Range("L94").Value = 10200
$ws_CUL.Range("M94").Value = -5144
$ws_CUL.Range("N94").Value = -11552
$ws_CUL.Range("H97").Value = 50001500
$ws_CUL.Range("I97").Value = 50001500
$ws_CUL.Range("K97").Value = 150004500
$ws_CUL.Range("M97").Value = -150004004
$ws_CUL.Range("H114").Value = 5895.3184
$ws_CUL.Range("I114").Value = 457.3
$ws_CUL.Range("J114").Value = 10427
$ws_CUL.Range("K114").Value = 1371.9
$ws_CUL.Range("L114").Value = 31281
$ws_CUL.Range("M114").Value = 1882.1
$ws_CUL.Range("N114").Value = -37789
$ws_CUL.Range("H122").Value = 2854.3171
$ws_CUL.Range("I122").Value = 420.5
$ws_CUL.Range("J122").Value = 4411.96
$ws_CUL.Range("K122").Value = 3784.5
$ws_CUL.Range("L122").Value = 39707.64
$ws_CUL.Range("M122").Value = -1334.5
$ws_CUL.Range("N122").Value = -44607.64
$ws_CUL.Range("H132").Value = 8717714
$ws_CUL.Range("I132").Value = 2875.8
$ws_CUL.Range("J132").Value = 13075133
$ws_CUL.Range("K132").Value = 25882.2
$ws_CUL.Range("L132").Value = 117676197
$ws_CUL.Range("M132").Value = -23352.2
$ws_CUL.Range("N132").Value = -117681257
$ws_GSM.Range("H70").Value = 5709.839
$ws_GSM.Range("I70").Value = 5852.3335
$ws_GSM.Range("J70").Value = 5221.2856
$ws_GSM.Range("K70").Value = 5852.3335
$ws_GSM.Range("L70").Value = 5221.2856
$ws_GSM.Range("M70").Value = -5582.3335
$ws_GSM.Range("N70").Value = -5761.2856
$ws_GSM.Range("H73").Value = 5709.839
$ws_GSM.Range("I73").Value = 5852.3335
$ws_GSM.Range("J73").Value = 5221.2856
$ws_GSM.Range("K73").Value = 5852.3335
$ws_GSM.Range("L73").Value = 5221.2856
$ws_GSM.Range("M73").Value = -4916.3335
$ws_GSM.Range("N73").Value = -7093.2856
$ws_LTW.Range("H136").Value = 6762.245
$ws_LTW.Range("I136").Value = 4300
$ws_LTW.Range("K136").Value = 12900
$ws_LTW.Range("M136").Value = -10350
$ws_WVR.Range("H81").Value = 1541.1818
$ws_WVR.Range("I81").Value = 1218.5714
$ws_WVR.Range("K81").Value = 2437.1428
$ws_WVR.Range("M81").Value = -1376.1428
$ws_WVR.Range("H84").Value = 1541.1818
$ws_WVR.Range("I84").Value = 1218.5714
$ws_WVR.Range("K84").Value = 12185.714
$ws_WVR.Range("M84").Value = -6881.714
$ws_WVR.Range("H103").Value = 29000
$ws_WVR.Range("J103").Value = 29000
$ws_WVR.Range("L103").Value = 29000
$ws_WVR.Range("N103").Value = -31344
$ws_WVR.Range("H132").Value = 14929.676
$ws_WVR.Range("I132").Value = 17032.13
$ws_WVR.Range("J132").Value = 2104.7
$ws_WVR.Range("K132").Value = 51096.39
$ws_WVR.Range("L132").Value = 6314.099999999999
$ws_WVR.Range("M132").Value = -48566.39
$ws_WVR.Range("N132").Value = -11374.1
$ws_WVR.Range("H136").Value = 10207267
$ws_WVR.Range("I136").Value = 4182.8
$ws_WVR.Range("J136").Value = 20835480
$ws_WVR.Range("K136").Value = 12548.4
$ws_WVR.Range("L136").Value = 62506440
$ws_WVR.Range("M136").Value = -9998.400000000001
$ws_WVR.Range("N136").Value = -62511540
$ws_ALC.Range("M138").Value = 3739.3
$ws_CUL.Range("N93").ClearContents()
